$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") moves forward one day (2024-12-04 -> 2024-12-05)
# for every data row still present after the edit (rows 2..34).
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 3).Value = 45631
}

# Row 35 (the last data row, "A 57397-2024") is removed entirely.
$ws.Rows.Item(35).Delete()

# The new last row (34) drops its explicit row height, matching the
# "last row has no custom height" pattern used throughout this sheet.
$ws.Rows.Item(34).AutoFit()
